$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8929305672645569
$ws.Range("B1").Value = 1.133575797080994
$ws.Range("C1").Value = 1.675683617591858
$ws.Range("D1").Value = 2.170413732528687
$ws.Range("E1").Value = 1.801521301269531
